$d = $word.ActiveDocument

# Find the paragraph that holds the "LOM3215 ... (Requisito)" line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOM3215: Física do Estado Sólido (Requisito)*") {
        $target = $p
        break
    }
}

# The three paragraphs that directly follow it are removed:
#   - a blank paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - the "© 2020 . Contact: ..." footer line
# The paragraph after those (another blank one, right before the page
# break) is left untouched.
$startPara = $target.Next()
$endPara = $startPara.Next().Next().Next()

$start = $startPara.Range.Start
$end = $endPara.Range.Start

$r = $d.Range($start, $end)
$r.Delete()
